# Applies the "Add files via upload" edit:
#  - Slide-number placeholder field preview text '<N deg>' -> '<#>' on the
#    slide master and every slide layout.
#  - Clears the authorship/contact text box (shape id 237) on slides
#    256, 259 and 261, leaving an empty paragraph (matches the author
#    stripping out the "Ashish Patel / Abonia Sojasingarayar / Updated:"
#    byline runs).
#  - Removes the stray LinkedIn / GitHub / Medium / HuggingFace logo
#    pictures that were left over on those same three slides.

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

function Set-SlideNumberFieldText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -eq "Slide Number" -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = [char]0x2039 + "#" + [char]0x203A
        }
    }
}

function Remove-ShapesByIds($shapes, $ids) {
    foreach ($id in $ids) {
        $sh = Get-ShapeById $shapes $id
        if ($sh -ne $null) {
            $sh.Delete()
        }
    }
}

$p = $ppt.ActivePresentation

# 1) Slide master + every layout: '<N deg>' -> '<#>' on the slide-number field.
$master = $p.SlideMaster
Set-SlideNumberFieldText $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-SlideNumberFieldText $layout.Shapes
}

# 2) Per-slide cleanup.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    $sldId = $s.SlideID

    # Clear the byline/contact text box (shape id 237) down to an empty
    # paragraph on the three affected slides.
    if ($sldId -eq 256 -or $sldId -eq 259 -or $sldId -eq 261) {
        $byline = Get-ShapeById $s.Shapes 237
        if ($byline -ne $null) {
            $byline.TextFrame.TextRange.Text = ""
        }
    }

    if ($sldId -eq 256) {
        Remove-ShapesByIds $s.Shapes @(1034, 1036, 23, 24, 72)
    } elseif ($sldId -eq 259) {
        Remove-ShapesByIds $s.Shapes @(10, 23, 24, 29, 1034, 1036)
    } elseif ($sldId -eq 261) {
        Remove-ShapesByIds $s.Shapes @(23, 24, 29, 1034, 1036)
    }
}
